# Commit: "Fruta / hortaliza, semanal"
# This weekly data refresh inserts one new Coliflor price observation
# (Primera + Segunda quality rows) into the Femacal de La Calera table.
# The new pair is inserted right before the existing row 542, pushing all
# subsequent rows down by two positions (dimension grows from R660 to R662).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 542-543; everything from old row 542 onward
# shifts down by two rows (old 542 -> 544, ..., old 660 -> 662).
$ws.Rows("542:543").Insert()

# New row 542 ("Primera" quality) for date 2022-05-30 (serial 44711).
$ws.Cells.Item(542, 1).Value = 3
$ws.Cells.Item(542, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(542, 3).Value = "Coquimbo"
$ws.Cells.Item(542, 4).Value = 44711
$ws.Cells.Item(542, 5).Value = 5
$ws.Cells.Item(542, 6).Value = 100112008
$ws.Cells.Item(542, 7).Value = "Coliflor"
$ws.Cells.Item(542, 8).Value = "Sin especificar"
$ws.Cells.Item(542, 9).Value = "Primera"
$ws.Cells.Item(542, 10).Value = 3100
$ws.Cells.Item(542, 11).Value = 850
$ws.Cells.Item(542, 12).Value = 900
$ws.Cells.Item(542, 13).Value = 874
$ws.Cells.Item(542, 14).Value = "`$/unidad"
$ws.Cells.Item(542, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(542, 16).Value = 874
$ws.Cells.Item(542, 17).Value = 1
$ws.Cells.Item(542, 18).Value = "Hortaliza"

# New row 543 ("Segunda" quality) for the same date 2022-05-30.
$ws.Cells.Item(543, 1).Value = 3
$ws.Cells.Item(543, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(543, 3).Value = "Coquimbo"
$ws.Cells.Item(543, 4).Value = 44711
$ws.Cells.Item(543, 5).Value = 5
$ws.Cells.Item(543, 6).Value = 100112008
$ws.Cells.Item(543, 7).Value = "Coliflor"
$ws.Cells.Item(543, 8).Value = "Sin especificar"
$ws.Cells.Item(543, 9).Value = "Segunda"
$ws.Cells.Item(543, 10).Value = 1200
$ws.Cells.Item(543, 11).Value = 650
$ws.Cells.Item(543, 12).Value = 650
$ws.Cells.Item(543, 13).Value = 650
$ws.Cells.Item(543, 14).Value = "`$/unidad"
$ws.Cells.Item(543, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(543, 16).Value = 650
$ws.Cells.Item(543, 17).Value = 1
$ws.Cells.Item(543, 18).Value = "Hortaliza"
